$d = $word.ActiveDocument

$d.Content.Find.Execute("74-7=", $true, $false, $false, $false, $false, $true, 1, $false, "90-20=", 2) | Out-Null
$d.Content.Find.Execute("5+46=", $true, $false, $false, $false, $false, $true, 1, $false, "47-17=", 2) | Out-Null
$d.Content.Find.Execute("35-35=", $true, $false, $false, $false, $false, $true, 1, $false, "23+56=", 2) | Out-Null
$d.Content.Find.Execute("85+6=", $true, $false, $false, $false, $false, $true, 1, $false, "22-8=", 2) | Out-Null
$d.Content.Find.Execute("74-37=", $true, $false, $false, $false, $false, $true, 1, $false, "46+6=", 2) | Out-Null
$d.Content.Find.Execute("75-25=", $true, $false, $false, $false, $false, $true, 1, $false, "52-35=", 2) | Out-Null
$d.Content.Find.Execute("9+0=", $true, $false, $false, $false, $false, $true, 1, $false, "35+54=", 2) | Out-Null
$d.Content.Find.Execute("5-0=", $true, $false, $false, $false, $false, $true, 1, $false, "19+44=", 2) | Out-Null
$d.Content.Find.Execute("16+62=", $true, $false, $false, $false, $false, $true, 1, $false, "81+17=", 2) | Out-Null
$d.Content.Find.Execute("71-13=", $true, $false, $false, $false, $false, $true, 1, $false, "75-2=", 2) | Out-Null
$d.Content.Find.Execute("37+12=", $true, $false, $false, $false, $false, $true, 1, $false, "8+86=", 2) | Out-Null
$d.Content.Find.Execute("20+23=", $true, $false, $false, $false, $false, $true, 1, $false, "23-10=", 2) | Out-Null
$d.Content.Find.Execute("13+56=", $true, $false, $false, $false, $false, $true, 1, $false, "74-50=", 2) | Out-Null
$d.Content.Find.Execute("86-28=", $true, $false, $false, $false, $false, $true, 1, $false, "59-28=", 2) | Out-Null
$d.Content.Find.Execute("12+0=", $true, $false, $false, $false, $false, $true, 1, $false, "15+66=", 2) | Out-Null
$d.Content.Find.Execute("77-75=", $true, $false, $false, $false, $false, $true, 1, $false, "96-89=", 2) | Out-Null
$d.Content.Find.Execute("81-22=", $true, $false, $false, $false, $false, $true, 1, $false, "78-5=", 2) | Out-Null
$d.Content.Find.Execute("76-53=", $true, $false, $false, $false, $false, $true, 1, $false, "9+73=", 2) | Out-Null
$d.Content.Find.Execute("72-61=", $true, $false, $false, $false, $false, $true, 1, $false, "12+26=", 2) | Out-Null
$d.Content.Find.Execute("32-13=", $true, $false, $false, $false, $false, $true, 1, $false, "51+15=", 2) | Out-Null
$d.Content.Find.Execute("78-11=", $true, $false, $false, $false, $false, $true, 1, $false, "19+45=", 2) | Out-Null
$d.Content.Find.Execute("44+12=", $true, $false, $false, $false, $false, $true, 1, $false, "6+51=", 2) | Out-Null
$d.Content.Find.Execute("33+27=", $true, $false, $false, $false, $false, $true, 1, $false, "83-59=", 2) | Out-Null
$d.Content.Find.Execute("96-9=", $true, $false, $false, $false, $false, $true, 1, $false, "41+55=", 2) | Out-Null
$d.Content.Find.Execute("9+41=", $true, $false, $false, $false, $false, $true, 1, $false, "51-47=", 2) | Out-Null
$d.Content.Find.Execute("13+19=", $true, $false, $false, $false, $false, $true, 1, $false, "55-30=", 2) | Out-Null
$d.Content.Find.Execute("20+37=", $true, $false, $false, $false, $false, $true, 1, $false, "75+12=", 2) | Out-Null
$d.Content.Find.Execute("47-41=", $true, $false, $false, $false, $false, $true, 1, $false, "91-17=", 2) | Out-Null
$d.Content.Find.Execute("88-32=", $true, $false, $false, $false, $false, $true, 1, $false, "80-12=", 2) | Out-Null
$d.Content.Find.Execute("3+9=", $true, $false, $false, $false, $false, $true, 1, $false, "71+8=", 2) | Out-Null
$d.Content.Find.Execute("12+8=", $true, $false, $false, $false, $false, $true, 1, $false, "75-56=", 2) | Out-Null
$d.Content.Find.Execute("42+50=", $true, $false, $false, $false, $false, $true, 1, $false, "79-43=", 2) | Out-Null
$d.Content.Find.Execute("96+2=", $true, $false, $false, $false, $false, $true, 1, $false, "48+9=", 2) | Out-Null
$d.Content.Find.Execute("59+20=", $true, $false, $false, $false, $false, $true, 1, $false, "5+39=", 2) | Out-Null
$d.Content.Find.Execute("6+55=", $true, $false, $false, $false, $false, $true, 1, $false, "92-37=", 2) | Out-Null
$d.Content.Find.Execute("82-24=", $true, $false, $false, $false, $false, $true, 1, $false, "19-13=", 2) | Out-Null
$d.Content.Find.Execute("92-44=", $true, $false, $false, $false, $false, $true, 1, $false, "30+2=", 2) | Out-Null
$d.Content.Find.Execute("97-16=", $true, $false, $false, $false, $false, $true, 1, $false, "1+89=", 2) | Out-Null
$d.Content.Find.Execute("6+66=", $true, $false, $false, $false, $false, $true, 1, $false, "65+9=", 2) | Out-Null
$d.Content.Find.Execute("32+29=", $true, $false, $false, $false, $false, $true, 1, $false, "53+16=", 2) | Out-Null
$d.Content.Find.Execute("40+45=", $true, $false, $false, $false, $false, $true, 1, $false, "14+67=", 2) | Out-Null
$d.Content.Find.Execute("78-35=", $true, $false, $false, $false, $false, $true, 1, $false, "0+52=", 2) | Out-Null
$d.Content.Find.Execute("49+46=", $true, $false, $false, $false, $false, $true, 1, $false, "51+2=", 2) | Out-Null
$d.Content.Find.Execute("88-51=", $true, $false, $false, $false, $false, $true, 1, $false, "34-25=", 2) | Out-Null
$d.Content.Find.Execute("97-73=", $true, $false, $false, $false, $false, $true, 1, $false, "3+29=", 2) | Out-Null
$d.Content.Find.Execute("78-27=", $true, $false, $false, $false, $false, $true, 1, $false, "49-25=", 2) | Out-Null
$d.Content.Find.Execute("21-17=", $true, $false, $false, $false, $false, $true, 1, $false, "73-17=", 2) | Out-Null
$d.Content.Find.Execute("18+10=", $true, $false, $false, $false, $false, $true, 1, $false, "68-18=", 2) | Out-Null
$d.Content.Find.Execute("76-63=", $true, $false, $false, $false, $false, $true, 1, $false, "64+6=", 2) | Out-Null
$d.Content.Find.Execute("30+57=", $true, $false, $false, $false, $false, $true, 1, $false, "71-39=", 2) | Out-Null
$d.Content.Find.Execute("39+17=", $true, $false, $false, $false, $false, $true, 1, $false, "37+46=", 2) | Out-Null
$d.Content.Find.Execute("29+69=", $true, $false, $false, $false, $false, $true, 1, $false, "1-0=", 2) | Out-Null
$d.Content.Find.Execute("58-25=", $true, $false, $false, $false, $false, $true, 1, $false, "93-45=", 2) | Out-Null
$d.Content.Find.Execute("23-13=", $true, $false, $false, $false, $false, $true, 1, $false, "50+17=", 2) | Out-Null
$d.Content.Find.Execute("18+59=", $true, $false, $false, $false, $false, $true, 1, $false, "50+1=", 2) | Out-Null
$d.Content.Find.Execute("61+15=", $true, $false, $false, $false, $false, $true, 1, $false, "20+3=", 2) | Out-Null
$d.Content.Find.Execute("86+11=", $true, $false, $false, $false, $false, $true, 1, $false, "26+11=", 2) | Out-Null
$d.Content.Find.Execute("23+31=", $true, $false, $false, $false, $false, $true, 1, $false, "11+13=", 2) | Out-Null
$d.Content.Find.Execute("55-39=", $true, $false, $false, $false, $false, $true, 1, $false, "2+36=", 2) | Out-Null
$d.Content.Find.Execute("5-1=", $true, $false, $false, $false, $false, $true, 1, $false, "3+23=", 2) | Out-Null
$d.Content.Find.Execute("65+26=", $true, $false, $false, $false, $false, $true, 1, $false, "70-1=", 2) | Out-Null
$d.Content.Find.Execute("17+80=", $true, $false, $false, $false, $false, $true, 1, $false, "5+54=", 2) | Out-Null
$d.Content.Find.Execute("58-55=", $true, $false, $false, $false, $false, $true, 1, $false, "64+15=", 2) | Out-Null
$d.Content.Find.Execute("11+6=", $true, $false, $false, $false, $false, $true, 1, $false, "46-14=", 2) | Out-Null
$d.Content.Find.Execute("60+7=", $true, $false, $false, $false, $false, $true, 1, $false, "45+37=", 2) | Out-Null
$d.Content.Find.Execute("57+21=", $true, $false, $false, $false, $false, $true, 1, $false, "18+25=", 2) | Out-Null
$d.Content.Find.Execute("80-1=", $true, $false, $false, $false, $false, $true, 1, $false, "73-72=", 2) | Out-Null
$d.Content.Find.Execute("47-46=", $true, $false, $false, $false, $false, $true, 1, $false, "26+72=", 2) | Out-Null
$d.Content.Find.Execute("93-72=", $true, $false, $false, $false, $false, $true, 1, $false, "1+5=", 2) | Out-Null
$d.Content.Find.Execute("12+32=", $true, $false, $false, $false, $false, $true, 1, $false, "67+1=", 2) | Out-Null
$d.Content.Find.Execute("31+14=", $true, $false, $false, $false, $false, $true, 1, $false, "29+33=", 2) | Out-Null
$d.Content.Find.Execute("27-0=", $true, $false, $false, $false, $false, $true, 1, $false, "69-46=", 2) | Out-Null
$d.Content.Find.Execute("82-78=", $true, $false, $false, $false, $false, $true, 1, $false, "31+5=", 2) | Out-Null
$d.Content.Find.Execute("90-53=", $true, $false, $false, $false, $false, $true, 1, $false, "70-36=", 2) | Out-Null
$d.Content.Find.Execute("96-42=", $true, $false, $false, $false, $false, $true, 1, $false, "55-30=", 2) | Out-Null
$d.Content.Find.Execute("4+24=", $true, $false, $false, $false, $false, $true, 1, $false, "89-55=", 2) | Out-Null
$d.Content.Find.Execute("57+16=", $true, $false, $false, $false, $false, $true, 1, $false, "92-66=", 2) | Out-Null
$d.Content.Find.Execute("32+11=", $true, $false, $false, $false, $false, $true, 1, $false, "38+22=", 2) | Out-Null
$d.Content.Find.Execute("64-15=", $true, $false, $false, $false, $false, $true, 1, $false, "25-13=", 2) | Out-Null
$d.Content.Find.Execute("91-48=", $true, $false, $false, $false, $false, $true, 1, $false, "22+49=", 2) | Out-Null
$d.Content.Find.Execute("82+13=", $true, $false, $false, $false, $false, $true, 1, $false, "33-28=", 2) | Out-Null
$d.Content.Find.Execute("87-82=", $true, $false, $false, $false, $false, $true, 1, $false, "96-35=", 2) | Out-Null
$d.Content.Find.Execute("82-21=", $true, $false, $false, $false, $false, $true, 1, $false, "92-45=", 2) | Out-Null
$d.Content.Find.Execute("10+27=", $true, $false, $false, $false, $false, $true, 1, $false, "7+43=", 2) | Out-Null
$d.Content.Find.Execute("51-39=", $true, $false, $false, $false, $false, $true, 1, $false, "7+49=", 2) | Out-Null
$d.Content.Find.Execute("95-33=", $true, $false, $false, $false, $false, $true, 1, $false, "97-17=", 2) | Out-Null
$d.Content.Find.Execute("20+73=", $true, $false, $false, $false, $false, $true, 1, $false, "21+3=", 2) | Out-Null
$d.Content.Find.Execute("29+13=", $true, $false, $false, $false, $false, $true, 1, $false, "78+3=", 2) | Out-Null
$d.Content.Find.Execute("8+48=", $true, $false, $false, $false, $false, $true, 1, $false, "13+67=", 2) | Out-Null
$d.Content.Find.Execute("47-44=", $true, $false, $false, $false, $false, $true, 1, $false, "52-14=", 2) | Out-Null
$d.Content.Find.Execute("36-35=", $true, $false, $false, $false, $false, $true, 1, $false, "50+20=", 2) | Out-Null
$d.Content.Find.Execute("76-23=", $true, $false, $false, $false, $false, $true, 1, $false, "72-8=", 2) | Out-Null
$d.Content.Find.Execute("8+37=", $true, $false, $false, $false, $false, $true, 1, $false, "80+19=", 2) | Out-Null
$d.Content.Find.Execute("83-39=", $true, $false, $false, $false, $false, $true, 1, $false, "32-8=", 2) | Out-Null
$d.Content.Find.Execute("11+63=", $true, $false, $false, $false, $false, $true, 1, $false, "24-2=", 2) | Out-Null
$d.Content.Find.Execute("6+60=", $true, $false, $false, $false, $false, $true, 1, $false, "8+33=", 2) | Out-Null
$d.Content.Find.Execute("96-61=", $true, $false, $false, $false, $false, $true, 1, $false, "35+11=", 2) | Out-Null
$d.Content.Find.Execute("32+60=", $true, $false, $false, $false, $false, $true, 1, $false, "90-24=", 2) | Out-Null
$d.Content.Find.Execute("97-87=", $true, $false, $false, $false, $false, $true, 1, $false, "61+33=", 2) | Out-Null
$d.Content.Find.Execute("14+65=", $true, $false, $false, $false, $false, $true, 1, $false, "45-28=", 2) | Out-Null
